# Overhaul / Version 3 style edit: add totals row to BMI table, add two
# more copies of the BMI lookup table (one plain, one with headers
# hidden), and switch which sheet/cell is active.

$wb = $excel.ActiveWorkbook
$wsPokemon = $wb.Worksheets.Item(1)
$wsBMI = $wb.Worksheets.Item(2)

# ---------------------------------------------------------------------
# 1. Pokémon sheet: it is no longer the active tab, selection moves to D4
# ---------------------------------------------------------------------
$wsPokemon.Activate()
$wsPokemon.Range("D4").Select()

# ---------------------------------------------------------------------
# 2. BMI sheet: add a totals row to the existing "BMI" table
# ---------------------------------------------------------------------
$loBMI = $wsBMI.ListObjects.Item(1)
$loBMI.ShowTotals = $true
$loBMI.TotalsRowRange.Cells.Item(1, 1).Value = "Total"
$loBMI.TotalsRowRange.Cells.Item(1, 2).Formula = "=SUBTOTAL(110,BMI[Max])"
$loBMI.TotalsRowRange.Cells.Item(1, 3).Formula = "=SUBTOTAL(103,BMI[Status])"

# ---------------------------------------------------------------------
# 3. BMI sheet: a second copy of the lookup table at A11:C17 (BMI_5)
# ---------------------------------------------------------------------
$wsBMI.Range("A11").Value = "Min"
$wsBMI.Range("B11").Value = "Max"
$wsBMI.Range("C11").Value = "Status"

$rows = @(
  @(0, 18.5, "Underweight"),
  @(18.5, 25, "Normal weight"),
  @(25, 30, "Overweight"),
  @(30, 35, "Obesity class 1"),
  @(35, 40, "Obesity class 2"),
  @(40, 10000, "Obesity class 3")
)

for ($i = 0; $i -lt 6; $i++) {
  $r = 12 + $i
  $wsBMI.Cells.Item($r, 1).Value = $rows[$i][0]
  $wsBMI.Cells.Item($r, 2).Value = $rows[$i][1]
  $wsBMI.Cells.Item($r, 3).Value = $rows[$i][2]
}

$loBMI5 = $wsBMI.ListObjects.Add(1, $wsBMI.Range("A11:C17"), $null, 1)
$loBMI5.Name = "BMI_5"

# ---------------------------------------------------------------------
# 4. BMI sheet: a label plus a third copy of the table, with its header
#    row hidden, at A21:C26 (BMI_6)
# ---------------------------------------------------------------------
$wsBMI.Range("A20").Value = "Headers are hidden"

# Seed row 21 with header-like text so the new table picks up the
# "Min" / "Max" / "Status" column names, then hide the header row and
# restore row 21 to being a normal data row.
$wsBMI.Cells.Item(21, 1).Value = "Min"
$wsBMI.Cells.Item(21, 2).Value = "Max"
$wsBMI.Cells.Item(21, 3).Value = "Status"

for ($i = 1; $i -lt 6; $i++) {
  $r = 21 + $i
  $wsBMI.Cells.Item($r, 1).Value = $rows[$i][0]
  $wsBMI.Cells.Item($r, 2).Value = $rows[$i][1]
  $wsBMI.Cells.Item($r, 3).Value = $rows[$i][2]
}

$loBMI6 = $wsBMI.ListObjects.Add(1, $wsBMI.Range("A21:C26"), $null, 1)
$loBMI6.Name = "BMI_6"
$loBMI6.ShowHeaders = $false
$loBMI6.Resize($wsBMI.Range("A21:C26"))

$wsBMI.Cells.Item(21, 1).Value = $rows[0][0]
$wsBMI.Cells.Item(21, 2).Value = $rows[0][1]
$wsBMI.Cells.Item(21, 3).Value = $rows[0][2]

# ---------------------------------------------------------------------
# 5. BMI sheet becomes the active tab, with A15 selected
# ---------------------------------------------------------------------
$wsBMI.Activate()
$wsBMI.Range("A15").Select()
